# Update "想去人数" (number of interested attendees) figures to reflect
# the newly generated gh-pages data output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions): rows 3 and 4 hold the two events whose
# attendee counts increased.
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 917
$wsExhibit.Range("F4").Value = 1671

# Sheet "全部类型" (All types): same two events appear again at rows 5 and 6.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value = 917
$wsAll.Range("F6").Value = 1671
